# Apply weekly update to "Hortaliza, Vega Modelo de Temuco - Arveja Verde" sheet
# Inserts two new daily-price rows (one near the top of the data block, one near
# the bottom) and shifts the existing rows down accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49; existing rows 49-85 shift down to 50-86.
$ws.Rows(49).Insert()

# Insert a second new row at 85 (old row 84 is now at 85, old row 85 is now at 86);
# this pushes them down to 86 and 87 respectively.
$ws.Rows(85).Insert()

# New row 49: Fecha 2022-01-06, 30 sacos @ 18000 avg, Precio $/Kg 720
$ws.Cells.Item(49,1).Value = 10
$ws.Cells.Item(49,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(49,3).Value = "La Araucanía"
$ws.Cells.Item(49,4).Value = 44567
$ws.Cells.Item(49,5).Value = 9
$ws.Cells.Item(49,6).Value = 100112022
$ws.Cells.Item(49,7).Value = "Arveja Verde"
$ws.Cells.Item(49,8).Value = "Sin especificar"
$ws.Cells.Item(49,9).Value = "Primera"
$ws.Cells.Item(49,10).Value = 30
$ws.Cells.Item(49,11).Value = 18000
$ws.Cells.Item(49,12).Value = 18000
$ws.Cells.Item(49,13).Value = 18000
$ws.Cells.Item(49,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(49,15).Value = "Región de La Araucanía"
$ws.Cells.Item(49,16).Value = 720
$ws.Cells.Item(49,17).Value = 25
$ws.Cells.Item(49,18).Value = "Hortaliza"

# New row 85: Fecha 2022-01-07, 50 sacos @ 18000 avg, Precio $/Kg 720
$ws.Cells.Item(85,1).Value = 10
$ws.Cells.Item(85,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(85,3).Value = "La Araucanía"
$ws.Cells.Item(85,4).Value = 44568
$ws.Cells.Item(85,5).Value = 9
$ws.Cells.Item(85,6).Value = 100112022
$ws.Cells.Item(85,7).Value = "Arveja Verde"
$ws.Cells.Item(85,8).Value = "Sin especificar"
$ws.Cells.Item(85,9).Value = "Primera"
$ws.Cells.Item(85,10).Value = 50
$ws.Cells.Item(85,11).Value = 18000
$ws.Cells.Item(85,12).Value = 18000
$ws.Cells.Item(85,13).Value = 18000
$ws.Cells.Item(85,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(85,15).Value = "Región de La Araucanía"
$ws.Cells.Item(85,16).Value = 720
$ws.Cells.Item(85,17).Value = 25
$ws.Cells.Item(85,18).Value = "Hortaliza"

